$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 52

$ws.Cells.Item($row, 1).Value = "DMXXFA"
$ws.Cells.Item($row, 2).Value = "Engranaje de rodillo de fusor superior para Samsung"
$ws.Cells.Item($row, 3).Value = "SCX 3200 3205 4016 4116 4200 4216 4300 4316 4321 4521 SF560, ML 1610 1640 1710 1740 1860"
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 120000
$ws.Cells.Item($row, 6).Value = 8
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Formula = "=(E52-D52)*G52"
$ws.Cells.Item($row, 9).Formula = "=D52*F52"
$ws.Cells.Item($row, 10).Value = 0
